$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 261; this shifts the existing rows 261-294
# down to 262-295, preserving all their data.
$ws.Rows.Item(261).Insert()

# Populate the newly inserted row 261 with the new record.
$ws.Cells.Item(261, 1).Value = 11
$ws.Cells.Item(261, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(261, 3).Value = "Bíobío"
$ws.Cells.Item(261, 4).Value = 45154
$ws.Cells.Item(261, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(261, 5).Value = 8
$ws.Cells.Item(261, 6).Value = "Fruta"
$ws.Cells.Item(261, 7).Value = 100108
$ws.Cells.Item(261, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(261, 9).Value = 100108005
$ws.Cells.Item(261, 10).Value = "Piña"
$ws.Cells.Item(261, 11).Value = "Caramelo"
$ws.Cells.Item(261, 12).Value = "Segunda"
$ws.Cells.Item(261, 13).Value = 150
$ws.Cells.Item(261, 14).Value = 22000
$ws.Cells.Item(261, 15).Value = 23000
$ws.Cells.Item(261, 16).Value = 22333
$ws.Cells.Item(261, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(261, 18).Value = "Ecuador"
$ws.Cells.Item(261, 19).Value = 1595
$ws.Cells.Item(261, 20).Value = 14
